$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert a new row at 38 as a duplicate of row 37 (the "last row" with special bottom border),
# shifting the old rows 38-43 down to 39-44.
$ws.Rows(38).Insert()
$ws.Range("B37:J37").Copy()
$ws.Range("B38:J38").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B37:J37").Copy()
$ws.Range("B38:J38").PasteSpecial(-4163)  # xlPasteValues

# Step 2: Demote row 37 back to "normal" row styling (copy formats from row 36, a normal data row)
$ws.Range("B36:J36").Copy()
$ws.Range("B37:J37").PasteSpecial(-4122)  # xlPasteFormats

# Step 3: Update new row 38's period value to "2509" (new period added)
$ws.Range("E38").Value = "2509"

# Step 4: Update totals - "VALOR MORA" total and "Cant. Periodos" count
$ws.Range("E11").Value = 1309620
$ws.Range("F13").Value = 6

# Step 5: Apply center horizontal alignment to the "Periodo Mora" column (E) for all data rows
$ws.Range("E16:E38").HorizontalAlignment = -4108  # xlCenter

Write-Host "Script completed"
